$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 7249.7
$ws.Range("I18").Value = 7555.222
$ws.Range("K18").Value = 7555.222
$ws.Range("M18").Value = -7271.222

$ws.Range("H55").Value = 82.333336
$ws.Range("I55").Value = 78.8
$ws.Range("K55").Value = 78.8
$ws.Range("M55").Value = 135.2

$ws.Range("H92").Value = 1672.1538
$ws.Range("I92").Value = 744.75
$ws.Range("J92").Value = 2084.3333
$ws.Range("K92").Value = 744.75
$ws.Range("L92").Value = 2084.3333
$ws.Range("M92").Value = 503.25
$ws.Range("N92").Value = -4580.3333

$ws.Range("H137").Value = 8489.678
$ws.Range("I137").Value = 9553.462
$ws.Range("K137").Value = 28660.386
$ws.Range("M137").Value = -26110.386

$ws.Range("H138").Value = 3506.6462
$ws.Range("I138").Value = 613
$ws.Range("J138").Value = 4096.093
$ws.Range("K138").Value = 1839
$ws.Range("L138").Value = 12288.279
$ws.Range("M138").Value = 3301
$ws.Range("N138").Value = -22568.279

$ws.Range("H141").Value = 4125.394
$ws.Range("I141").Value = 3272.8572
$ws.Range("K141").Value = 9818.571599999999
$ws.Range("M141").Value = -4638.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2367.426
$ws.Range("I32").Value = 2356.84
$ws.Range("K32").Value = 2356.84
$ws.Range("M32").Value = -2069.84

$ws.Range("H45").Value = 90907.30499999999
$ws.Range("I45").Value = 128102.5
$ws.Range("J45").Value = 5889.7144
$ws.Range("K45").Value = 128102.5
$ws.Range("L45").Value = 5889.7144
$ws.Range("M45").Value = -127725.5
$ws.Range("N45").Value = -6643.7144

$ws.Range("H61").Value = 3259.5588
$ws.Range("I61").Value = 3150.3438
$ws.Range("K61").Value = 3150.3438
$ws.Range("M61").Value = -2938.3438

$ws.Range("H74").Value = 3979.5334
$ws.Range("I74").Value = 2831.0264
$ws.Range("J74").Value = 10214.286
$ws.Range("K74").Value = 2831.0264
$ws.Range("L74").Value = 10214.286
$ws.Range("M74").Value = -1957.0264
$ws.Range("N74").Value = -11962.286

$ws.Range("H77").Value = 3979.5334
$ws.Range("I77").Value = 2831.0264
$ws.Range("J77").Value = 10214.286
$ws.Range("K77").Value = 14155.132
$ws.Range("L77").Value = 51071.43
$ws.Range("M77").Value = -9787.132000000001
$ws.Range("N77").Value = -59807.43

$ws.Range("H102").Value = 23475.8
$ws.Range("I102").Value = 35425.332
$ws.Range("J102").Value = 5551.5
$ws.Range("K102").Value = 35425.332
$ws.Range("L102").Value = 5551.5
$ws.Range("M102").Value = -33803.332
$ws.Range("N102").Value = -8795.5

$ws.Range("H110").Value = 2446
$ws.Range("I110").Value = 1687.4166
$ws.Range("J110").Value = 4721.75
$ws.Range("K110").Value = 1687.4166
$ws.Range("L110").Value = 4721.75
$ws.Range("M110").Value = 357.5834
$ws.Range("N110").Value = -8811.75

$ws.Range("H122").Value = 705799.9399999999
$ws.Range("I122").Value = 5500
$ws.Range("J122").Value = 2806699.8
$ws.Range("K122").Value = 16500
$ws.Range("L122").Value = 8420099.399999999
$ws.Range("M122").Value = -14050
$ws.Range("N122").Value = -8424999.399999999

$ws.Range("H123").Value = 117777.5
$ws.Range("J123").Value = 117777.5
$ws.Range("L123").Value = 117777.5
$ws.Range("N123").Value = -127577.5

$ws.Range("H132").Value = 1678.4773
$ws.Range("I132").Value = 963.1142599999999
$ws.Range("J132").Value = 4460.4443
$ws.Range("K132").Value = 2889.34278
$ws.Range("L132").Value = 13381.3329
$ws.Range("M132").Value = -359.3427799999999
$ws.Range("N132").Value = -18441.3329

$ws.Range("H133").Value = 77499.5
$ws.Range("I133").Value = 76999
$ws.Range("K133").Value = 76999
$ws.Range("M133").Value = -74469

$ws.Range("H136").Value = 3259.5588
$ws.Range("I136").Value = 3150.3438
$ws.Range("K136").Value = 9451.0314
$ws.Range("M136").Value = -6901.0314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 95154.09
$ws.Range("I105").Value = 144671
$ws.Range("K105").Value = 144671
$ws.Range("M105").Value = -142924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2090.5908
$ws.Range("I58").Value = 1648.1818
$ws.Range("J58").Value = 2533
$ws.Range("K58").Value = 1648.1818
$ws.Range("L58").Value = 2533
$ws.Range("M58").Value = -1445.1818
$ws.Range("N58").Value = -2939

$ws.Range("H111").Value = 65000
$ws.Range("J111").Value = 65000
$ws.Range("L111").Value = 65000
$ws.Range("N111").Value = -73180

$ws.Range("H134").Value = 3422.689
$ws.Range("I134").Value = 1802.9143
$ws.Range("K134").Value = 5408.742899999999
$ws.Range("M134").Value = -2873.742899999999

$ws.Range("H136").Value = 2090.5908
$ws.Range("I136").Value = 1648.1818
$ws.Range("J136").Value = 2533
$ws.Range("K136").Value = 4944.5454
$ws.Range("L136").Value = 7599
$ws.Range("M136").Value = -2394.5454
$ws.Range("N136").Value = -12699

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 3622.7368
$ws.Range("I41").Value = 6183.2
$ws.Range("J41").Value = 777.7778
$ws.Range("K41").Value = 18549.6
$ws.Range("L41").Value = 2333.3334
$ws.Range("M41").Value = -18211.6
$ws.Range("N41").Value = -3009.3334

$ws.Range("H44").Value = 12793.833
$ws.Range("I44").Value = 509.33334
$ws.Range("J44").Value = 16888.666
$ws.Range("K44").Value = 1528.00002
$ws.Range("L44").Value = 50665.99800000001
$ws.Range("M44").Value = -1130.00002
$ws.Range("N44").Value = -51461.99800000001

$ws.Range("H51").Value = 2390.1428
$ws.Range("J51").Value = 3299.5
$ws.Range("L51").Value = 9898.5
$ws.Range("N51").Value = -10818.5

$ws.Range("H113").Value = 1335.9231
$ws.Range("J113").Value = 1552.5714
$ws.Range("L113").Value = 4657.7142
$ws.Range("N113").Value = -8997.7142

$ws.Range("H131").Value = 83335256
$ws.Range("I131").Value = 250000900
$ws.Range("K131").Value = 750002700
$ws.Range("M131").Value = -749997660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11184
$ws.Range("I80").Value = 12620.8
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 12620.8
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -11622.8
$ws.Range("N80").Value = -5996

$ws.Range("H83").Value = 11184
$ws.Range("I83").Value = 12620.8
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 63104
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -58112
$ws.Range("N83").Value = -29984

$ws.Range("H97").Value = 7555.4414
$ws.Range("I97").Value = 10281.046
$ws.Range("J97").Value = 2558.5
$ws.Range("K97").Value = 10281.046
$ws.Range("L97").Value = 2558.5
$ws.Range("M97").Value = -9785.046
$ws.Range("N97").Value = -3550.5

$ws.Range("H102").Value = 7454.44
$ws.Range("I102").Value = 8729.526
$ws.Range("K102").Value = 8729.526
$ws.Range("M102").Value = -7107.526

$ws.Range("H126").Value = 26844.846
$ws.Range("I126").Value = 55996.332
$ws.Range("J126").Value = 18099.4
$ws.Range("K126").Value = 167988.996
$ws.Range("L126").Value = 54298.2
$ws.Range("M126").Value = -165518.996
$ws.Range("N126").Value = -59238.2

$ws.Range("H132").Value = 3039.675
$ws.Range("I132").Value = 2851.1936
$ws.Range("J132").Value = 3688.889
$ws.Range("K132").Value = 8553.5808
$ws.Range("L132").Value = 11066.667
$ws.Range("M132").Value = -6023.5808
$ws.Range("N132").Value = -16126.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26583.578
$ws.Range("J7").Value = 6832.6665
$ws.Range("L7").Value = 6832.6665
$ws.Range("N7").Value = -7056.6665

$ws.Range("H16").Value = 3374.8333
$ws.Range("I16").Value = 4241.2856
$ws.Range("J16").Value = 342.25
$ws.Range("K16").Value = 4241.2856
$ws.Range("L16").Value = 342.25
$ws.Range("M16").Value = -4071.2856
$ws.Range("N16").Value = -682.25

$ws.Range("H40").Value = 36166.133
$ws.Range("I40").Value = 59732.715
$ws.Range("J40").Value = 15545.375
$ws.Range("K40").Value = 59732.715
$ws.Range("L40").Value = 15545.375
$ws.Range("M40").Value = -59596.715
$ws.Range("N40").Value = -15817.375

$ws.Range("H55").Value = 1741.3636
$ws.Range("I55").Value = 290.2857
$ws.Range("J55").Value = 4280.75
$ws.Range("K55").Value = 290.2857
$ws.Range("L55").Value = 4280.75
$ws.Range("M55").Value = -117.2857
$ws.Range("N55").Value = -4626.75

$ws.Range("H93").Value = 5136.5
$ws.Range("I93").Value = 6816.353
$ws.Range("J93").Value = 1056.8572
$ws.Range("K93").Value = 6816.353
$ws.Range("L93").Value = 1056.8572
$ws.Range("M93").Value = -5568.353
$ws.Range("N93").Value = -3552.8572

$ws.Range("H126").Value = 26583.578
$ws.Range("J126").Value = 6832.6665
$ws.Range("L126").Value = 20497.9995
$ws.Range("N126").Value = -25437.9995

$ws.Range("H132").Value = 517003.8
$ws.Range("I132").Value = 785164.0600000001
$ws.Range("J132").Value = 7499.4
$ws.Range("K132").Value = 2355492.18
$ws.Range("L132").Value = 22498.2
$ws.Range("M132").Value = -2352962.18
$ws.Range("N132").Value = -27558.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3100
$ws.Range("I14").Value = 5000
$ws.Range("J14").Value = 2466.6667
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 2466.6667
$ws.Range("M14").Value = -4832
$ws.Range("N14").Value = -2802.6667

$ws.Range("H107").Value = 13621.208
$ws.Range("I107").Value = 1518.25
$ws.Range("J107").Value = 37827.125
$ws.Range("K107").Value = 4554.75
$ws.Range("L107").Value = 113481.375
$ws.Range("M107").Value = -2634.75
$ws.Range("N107").Value = -117321.375

$ws.Range("H122").Value = 5455.7144
$ws.Range("I122").Value = 4188.05
$ws.Range("K122").Value = 12564.15
$ws.Range("M122").Value = -10114.15

$ws.Range("H126").Value = 37625.25
$ws.Range("I126").Value = 51439.125
$ws.Range("J126").Value = 9997.5
$ws.Range("K126").Value = 154317.375
$ws.Range("L126").Value = 29992.5
$ws.Range("M126").Value = -151847.375
$ws.Range("N126").Value = -34932.5

$ws.Range("H132").Value = 10841.818
$ws.Range("I132").Value = 11556.527
$ws.Range("K132").Value = 34669.581
$ws.Range("M132").Value = -32139.581
